$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data scraped by the GitHub Actions job
# D-column values that look numeric must be forced to Text format first so
# Excel does not silently coerce/round them (e.g. "5.360" -> 5.36).
$ws.Range("D2").Value = "30.203.74"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.927.74"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.94"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7116"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3201"
$ws.Range("E8").Value = "  -2.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.35"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07081"
$ws.Range("E10").Value = "  +3.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7893"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07918"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "1.925.98"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.360"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.86"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "30.225.57"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "255.15"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008004"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.750"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "2.183.31"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.799"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.515"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.75"
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.00"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.266"
$ws.Range("E28").Value = "  -5.18%  "
$ws.Range("E29").Value = "  -5.80%  "
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.526"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.380"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.105"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05146"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.264"
$ws.Range("E35").Value = "  +3.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7427"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.763"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01946"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "77.44"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.322"
$ws.Range("E41").Value = "  -4.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4472"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.972"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8419"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.61"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.678"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.416"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.29"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06116"
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.863"
$ws.Range("E51").Value = "  +8.29%  "
